# Updated cryptos list on Fri Dec 22 05:10:14 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format Price (column D) cells whose new values would otherwise be
# auto-parsed as numbers, so they stay plain text like the rest of column D.
$numericPriceCells = @("D5", "D6", "D9", "D10", "D12", "D16", "D19", "D21", "D22", "D23", "D24", "D26", "D27", "D28", "D29", "D31", "D32", "D33", "D34", "D36", "D37", "D38", "D40", "D41", "D42", "D43", "D44", "D46", "D47", "D50")
foreach ($addr in $numericPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "44.184.78"
$ws.Range("E2").Value = "  +1.26%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.259.15"
$ws.Range("E3").Value = "  +2.69%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.05%  "

# Row 5 - Solana
$ws.Range("D5").Value = "98.56"
$ws.Range("E5").Value = "  +16.15%  "

# Row 6 - BNB
$ws.Range("D6").Value = "273.92"
$ws.Range("E6").Value = "  +6.43%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +1.02%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.02%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "0.630"
$ws.Range("E9").Value = "  +5.84%  "

# Row 10 - Avalanche
$ws.Range("D10").Value = "48.35"
$ws.Range("E10").Value = "  +7.09%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +2.82%  "

# Row 12 - Polkadot
$ws.Range("D12").Value = "8.18"
$ws.Range("E12").Value = "  +13.63%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +0.70%  "

# Row 14 - Chainlink
$ws.Range("E14").Value = "  +8.36%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "2.595.11"
$ws.Range("E15").Value = "  +2.53%  "

# Row 16 - Polygon
$ws.Range("D16").Value = "0.839"
$ws.Range("E16").Value = "  +6.87%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.262.49"
$ws.Range("E17").Value = "  +2.91%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "44.134.88"
$ws.Range("E18").Value = "  +1.22%  "

# Row 19 - ShibaInu
$ws.Range("D19").Value = "0.0000107"
$ws.Range("E19").Value = "  +4.02%  "

# Row 20 - Uniswap
$ws.Range("E20").Value = "  +5.72%  "

# Row 21 - Litecoin
$ws.Range("D21").Value = "71.01"
$ws.Range("E21").Value = "  +1.79%  "

# Row 22 - ImmutableX/InternetComputer swap
$ws.Range("B22").Value = "InternetComputer(DFINITY)"
$ws.Range("C22").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D22").Value = "10.56"
$ws.Range("E22").Value = "  +16.39%  "

# Row 23 - InternetComputer/ImmutableX swap
$ws.Range("B23").Value = "ImmutableX"
$ws.Range("C23").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D23").Value = "2.38"
$ws.Range("E23").Value = "  +0.18%  "

# Row 24 - BitcoinCash
$ws.Range("D24").Value = "235.29"
$ws.Range("E24").Value = "  +1.61%  "

# Row 25 - Dai
$ws.Range("E25").Value = "  +0.02%  "

# Row 26 - Cosmos
$ws.Range("D26").Value = "11.46"
$ws.Range("E26").Value = "  +7.72%  "

# Row 27 - PancakeSwap
$ws.Range("D27").Value = "2.51"
$ws.Range("E27").Value = "  +13.24%  "

# Row 28 - WEMIXToken
$ws.Range("D28").Value = "3.49"
$ws.Range("E28").Value = "  -1.97%  "

# Row 29 - InjectiveProtocol
$ws.Range("D29").Value = "40.39"
$ws.Range("E29").Value = "  +3.94%  "

# Row 30 - Toncoin
$ws.Range("E30").Value = "  +0.41%  "

# Row 31 - Monero
$ws.Range("D31").Value = "173.84"
$ws.Range("E31").Value = "  -0.06%  "

# Row 32 - Hedera
$ws.Range("D32").Value = "0.0920"
$ws.Range("E32").Value = "  +7.00%  "

# Row 33 - EthereumClassic
$ws.Range("D33").Value = "21.21"
$ws.Range("E33").Value = "  +4.10%  "

# Row 34 - Filecoin
$ws.Range("D34").Value = "5.75"
$ws.Range("E34").Value = "  +8.01%  "

# Row 35 - Stellar
$ws.Range("E35").Value = "  +1.82%  "

# Row 36 - Kaspa
$ws.Range("D36").Value = "0.115"
$ws.Range("E36").Value = "  +3.31%  "

# Row 37 - VeChain
$ws.Range("D37").Value = "0.0357"
$ws.Range("E37").Value = "  -0.71%  "

# Row 38 - RenderToken
$ws.Range("D38").Value = "4.41"
$ws.Range("E38").Value = "  -1.73%  "

# Row 39 - NEARProtocol
$ws.Range("E39").Value = "  +24.10%  "

# Row 40 - Algorand
$ws.Range("D40").Value = "0.255"
$ws.Range("E40").Value = "  +28.14%  "

# Row 41 - Celestia
$ws.Range("D41").Value = "12.65"
$ws.Range("E41").Value = "  +1.89%  "

# Row 42 - LidoDAOToken
$ws.Range("D42").Value = "2.19"
$ws.Range("E42").Value = "  +4.45%  "

# Row 43 - MultiversX
$ws.Range("D43").Value = "62.58"

# Row 44 - THORChain
$ws.Range("D44").Value = "5.48"
$ws.Range("E44").Value = "  +0.25%  "

# Row 45 - Cronos
$ws.Range("E45").Value = "  +6.01%  "

# Row 46 - FraxShare
$ws.Range("D46").Value = "8.58"
$ws.Range("E46").Value = "  +3.27%  "

# Row 47 - Aave
$ws.Range("D47").Value = "100.69"
$ws.Range("E47").Value = "  +0.62%  "

# Row 48 - ARBITRUM
$ws.Range("E48").Value = "  +4.60%  "

# Row 49 - TrustWalletToken
$ws.Range("E49").Value = "  +0.93%  "

# Row 50 - WOONetwork
$ws.Range("D50").Value = "0.433"
$ws.Range("E50").Value = "  +0.45%  "

# Row 51 - RocketPoolETH
$ws.Range("D51").Value = "2.478.21"
$ws.Range("E51").Value = "  +2.53%  "
